$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 5: header-style row for the new "send http get request" test case ---
$ws.Range("B5").Value = "url"
$ws.Range("C5").Value = "expected response"
$ws.Range("D5").Value = "json"
$ws.Range("E5").Value = "expected response regex"

# --- New row 6: the actual test data for response / regex assertions ---
$ws.Range("A6").Value = "send http post request"
$ws.Range("B6").Value = "http://requestb.in/re7qq8re"
$ws.Range("C6").Value = "ok"
$ws.Range("D6").Value = '{"query":"data"}'
$ws.Range("E6").Value = "o."

# B6 (the url) gets the blue "hyperlink-ish" Arial font used elsewhere for urls
$ws.Range("B6").Font.Color = 16711680
$ws.Range("B6").Font.Name = "Arial"
$ws.Range("B6").Font.Size = 10

# --- Column widths: columns A-C got narrower, D/E are new, F.. stays default ---
$ws.Columns.Item(1).ColumnWidth = 19.7908163265306
$ws.Columns.Item(2).ColumnWidth = 15.4336734693878
$ws.Columns.Item(3).ColumnWidth = 43.4744897959184
$ws.Columns.Item(4).ColumnWidth = 13.7857142857143
$ws.Columns.Item(5).ColumnWidth = 21.6785714285714

# --- View: selection moved to the newly added F5:F6 range ---
$ws.Range("F5:F6").Select()
$excel.Goto($ws.Range("F6"), $false)

Write-Output "edit applied"
